# Updated symbol list on Tue Dec 27 21:36:20 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are stored as text in this sheet (not numbers), so a
# leading apostrophe forces Excel to keep them as literal text instead of coercing
# them to floating-point numbers (which would also mangle trailing zeros).

$ws.Range("D2").Value = "'245.31"
$ws.Range("D3").Value = "'24.01"
$ws.Range("D5").Value = "'0.05818"
$ws.Range("D6").Value = "'6.467"
$ws.Range("D7").Value = "'3.360"
$ws.Range("D8").Value = "'0.8096"
$ws.Range("D9").Value = "'0.9181"
$ws.Range("D10").Value = "'0.1405"
$ws.Range("D11").Value = "'0.07349"
$ws.Range("D12").Value = "'0.03177"
$ws.Range("D13").Value = "'0.03077"
$ws.Range("D14").Value = "'0.09371"
$ws.Range("D15").Value = "'3.849"
$ws.Range("D16").Value = "'0.001557"
$ws.Range("D17").Value = "'0.04701"
$ws.Range("D18").Value = "'0.0005980"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.006123"
$ws.Range("D20").Value = "'0.001245"
$ws.Range("D21").Value = "'0.004689"
$ws.Range("D23").Value = "'3.592"
$ws.Range("D26").Value = "'0.1327"
$ws.Range("D40").Value = "'0.03841"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1065"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002750"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003066"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
$ws.Range("D44").Value = "'0.009038"
$ws.Range("D45").Value = "'0.00005250"
$ws.Range("D48").Value = "'0.001837"
